$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price/volume data (and a few swapped coin name/link pairs)
# to reflect the latest scrape. Price column (D) values must stay as literal
# text (e.g. "0.320", "23.20", "0.0000276") rather than being auto-converted
# to numbers by Excel, so we force text format before assigning, then restore
# the default "Normal" style so no stray formatting is introduced.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "84.431.53"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +4.56%  "
# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.290.08"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.98%  "
# Row 4
$ws.Range("E4").Value = "  +0.08%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.81"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.30%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "634.96"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.08%  "
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.320"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +17.20%  "
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.591"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.80%  "
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "3.291.04"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.09%  "
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.595"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.50%  "
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000276"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.85%  "
# Row 13
$ws.Range("E13").Value = "  -0.14%  "
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.894.20"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.15%  "
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "33.99"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.00%  "
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.41"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.57%  "
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "84.480.00"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.95%  "
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.291.97"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.27%  "
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.56"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.17%  "
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.16"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.05%  "
# Row 21
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "439.21"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.22%  "
# Row 22
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.18"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.64%  "
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.23"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.14%  "
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.37"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.54%  "
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.42"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +10.73%  "
# Row 26
$ws.Range("E26").Value = "  +10.93%  "
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.453.87"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.12%  "
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "77.81"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.29%  "
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0000132"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.37%  "
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.996"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.26%  "
# Row 31
$ws.Range("B31").Value = "Cronos"
$ws.Range("C31").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.167"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +34.66%  "
# Row 32
$ws.Range("B32").Value = "Bittensor"
$ws.Range("C32").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "604.82"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +6.70%  "
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "9.23"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.59%  "
# Row 34
$ws.Range("B34").Value = "Binance-PegBSC-USD"
$ws.Range("C34").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.998"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.09%  "
# Row 35
$ws.Range("B35").Value = "Fetch.AI"
$ws.Range("C35").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.57"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.87%  "
# Row 36
$ws.Range("E36").Value = "  -0.76%  "
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.151"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.11%  "
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "23.20"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.59%  "
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.34"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +7.86%  "
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.997"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.16%  "
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.416"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.37%  "
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.13"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +12.50%  "
# Row 43
$ws.Range("B43").Value = "WhiteBITCoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "20.92"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.89%  "
# Row 44
$ws.Range("B44").Value = "Stacks"
$ws.Range("C44").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.03"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +9.41%  "
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "159.66"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.01%  "
# Row 46
$ws.Range("E46").Value = "  +0.05%  "
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "189.32"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.93%  "
# Row 48
$ws.Range("B48").Value = "ImmutableX"
$ws.Range("C48").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.35"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.00%  "
# Row 49
$ws.Range("B49").Value = "OKB"
$ws.Range("C49").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "44.77"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.69%  "
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.787"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.44%  "
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "26.60"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.78%  "
